# Daily update at 8 AM UTC
# Appends the next day's row of data (row 91) to the "Wins Over Time" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

$ws.Cells.Item($row, 1).Value = 46040
$ws.Cells.Item($row, 2).Value = 210
$ws.Cells.Item($row, 3).Value = 218
$ws.Cells.Item($row, 4).Value = 204

# Match the date-style formatting used by the rest of column A.
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
